$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 229
$ws.Range("I2").Value = 554
$ws.Range("J2").Value = 2497
$ws.Range("K2").Value = 21
$ws.Range("L2").Value = 681
$ws.Range("N2").Value = 392
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 7
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 28
$ws.Range("S2").Value = 270
$ws.Range("T2").Value = 452
$ws.Range("U2").Value = 38
$ws.Range("V2").Value = 3835
$ws.Range("W2").Value = 4
$ws.Range("X2").Value = 3732
$ws.Range("Y2").Value = 6
$ws.Range("Z2").Value = 53
$ws.Range("AA2").Value = 31
